$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing 3 data rows -----------------------------------
$ws.Range("A1").Value = "Get patient's details: patient requests patient details"
$ws.Range("B1").Value = "BE"
$ws.Range("C1").Value = "GRD-713"

$ws.Range("A2").Value = "Carer confirmation: patient adds to his list not existing in DB carer"
$ws.Range("B2").Value = "BE"
$ws.Range("C2").Value = "GRD-713"

$ws.Range("A3").Value = "Carer confirmation: carer added new patient, patient confirmes carer"
$ws.Range("B3").Value = "BE"
$ws.Range("C3").Value = "GRD-713"

# --- Bring row 1's "no border / wrap+top" look onto the two new rows --
# before populating them, so the new cells inherit the same cell style as
# the existing column-A cells (minus the thin border, removed below).
$ws.Range("A1").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New rows 4 and 5 ---------------------------------------------------
$ws.Range("A4").Value = "User activation: inactive patient log in to app"
$ws.Range("B4").Value = "BE"
$ws.Range("C4").Value = "GRD-713"

$ws.Range("A5").Value = "User activation: inactive carer log in to app"
$ws.Range("B5").Value = "BE"
$ws.Range("C5").Value = "GRD-713"

# --- Column A: keep wrap + top alignment, drop the thin box border -----
$colA = $ws.Range("A1:A5")
$colA.WrapText = $true
$colA.VerticalAlignment = -4160
$colA.Borders.LineStyle = -4142

# --- Row heights: clear the old fixed 30pt rows, go back to automatic --
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()

# --- Column A width ------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 106.67

# --- Selection: whole first row (header row) ----------------------------
$ws.Range("A1:XFD1").Select()
